$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($addr, $val)
    # Force the assignment to remain a text value (matching the workbook's
    # original inlineStr/text cell type) even when the new text happens to
    # look like a plain number (e.g. "212.68"). NumberFormat="@" makes Excel
    # store the literal text instead of coercing it to a Number; ClearFormats()
    # afterwards drops the now-unneeded text format so the cell's style stays
    # at its original default (no stray formatting is introduced).
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.ClearFormats()
}

# Row 2 - Bitcoin
Set-TextValue "D2" "27.570.05"
$ws.Range("E2").Value = "  -0.08%  "

# Row 3 - Ethereum
Set-TextValue "D3" "1.647.71"
$ws.Range("E3").Value = "  -0.98%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  -0.10%  "

# Row 5 - BNB
Set-TextValue "D5" "212.68"
$ws.Range("E5").Value = "  -1.28%  "

# Row 6 - XRP
Set-TextValue "D6" "0.532"
$ws.Range("E6").Value = "  +4.59%  "

# Row 7 - USDC
$ws.Range("E7").Value = "  -0.07%  "

# Row 8 - Solana
Set-TextValue "D8" "23.65"
$ws.Range("E8").Value = "  -2.93%  "

# Row 10 - Dogecoin
$ws.Range("E10").Value = "  -1.31%  "

# Row 11 - TRON
Set-TextValue "D11" "0.0892"
$ws.Range("E11").Value = "  +1.53%  "

# Row 12 - WrappedliquidstakedEther2.0
Set-TextValue "D12" "1.880.71"
$ws.Range("E12").Value = "  -1.01%  "

# Row 13 - WrappedEther
Set-TextValue "D13" "1.659.62"
$ws.Range("E13").Value = "  -0.24%  "

# Row 14 - Polygon
Set-TextValue "D14" "0.591"
$ws.Range("E14").Value = "  +3.94%  "

# Row 15 - Polkadot
$ws.Range("E15").Value = "  -2.18%  "

# Row 16 - Litecoin
$ws.Range("E16").Value = "  -2.52%  "

# Row 17 - WrappedBTC
Set-TextValue "D17" "27.526.47"
$ws.Range("E17").Value = "  -0.18%  "

# Row 18 - BitcoinCash
Set-TextValue "D18" "232.29"
$ws.Range("E18").Value = "  -3.61%  "

# Row 19 - ShibaInu
$ws.Range("D19").Value = "0.0₃0726"
$ws.Range("E19").Value = "  -0.76%  "

# Row 20 - Chainlink
Set-TextValue "D20" "7.57"
$ws.Range("E20").Value = "  -1.07%  "

# Row 21 - Dai
$ws.Range("E21").Value = "  -0.13%  "

# Row 22 - Uniswap
$ws.Range("E22").Value = "  -3.67%  "

# Row 23 - Avalanche
Set-TextValue "D23" "9.81"
$ws.Range("E23").Value = "  +4.45%  "

# Row 24 - Toncoin
$ws.Range("E24").Value = "  -1.76%  "

# Row 25 - Monero
Set-TextValue "D25" "148.71"

# Row 26 - Cosmos
Set-TextValue "D26" "7.04"
$ws.Range("E26").Value = "  -3.00%  "

# Row 27 - Stellar
$ws.Range("E27").Value = "  +1.56%  "

# Row 28 - BinanceUSD
$ws.Range("E28").Value = "  -0.19%  "

# Row 29 - EthereumClassic
$ws.Range("E29").Value = "  -4.49%  "

# Row 30 - PancakeSwap
$ws.Range("E30").Value = "  -2.75%  "

# Row 31 - Hedera
$ws.Range("E31").Value = "  -3.43%  "

# Row 32 - Filecoin
Set-TextValue "D32" "3.33"
$ws.Range("E32").Value = "  -0.81%  "

# Row 33 - InternetComputer(DFINITY)
$ws.Range("E33").Value = "  +1.04%  "

# Row 34 - Maker
Set-TextValue "D34" "1.424.99"
$ws.Range("E34").Value = "  -2.44%  "

# Row 35 - LidoDAOToken
$ws.Range("E35").Value = "  +1.17%  "

# Row 36 - HuobiToken
$ws.Range("E36").Value = "  +0.16%  "

# Row 37 - ImmutableX
Set-TextValue "D37" "0.570"
$ws.Range("E37").Value = "  -0.56%  "

# Row 38 - ARBITRUM
$ws.Range("E38").Value = "  -4.54%  "

# Row 39 - VeChain
$ws.Range("E39").Value = "  -3.19%  "

# Row 40 - WEMIXToken
$ws.Range("E40").Value = "  -0.16%  "

# Row 41 - PaxDollar
$ws.Range("E41").Value = "  -0.05%  "

# Row 42 - TrustWalletToken
Set-TextValue "D42" "0.820"
$ws.Range("E42").Value = "  +3.48%  "

# Row 43 - now FraxShare (was mCoin)
$ws.Range("B43").Value = "FraxShare"
$ws.Range("C43").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
Set-TextValue "D43" "5.54"
$ws.Range("E43").Value = "  +2.20%  "

# Row 44 - now mCoin (was FraxShare)
$ws.Range("B44").Value = "mCoin"
$ws.Range("C44").Value = "https://coinranking.com/coin/fzVgyjBcRc9+mcoin-mcoin"
Set-TextValue "D44" "2.46"
$ws.Range("E44").Value = "  -3.31%  "

# Row 45 - MXToken
$ws.Range("E45").Value = "  +1.04%  "

# Row 46 - Aave
Set-TextValue "D46" "65.12"
$ws.Range("E46").Value = "  -6.96%  "

# Row 47 - RocketPoolETH
Set-TextValue "D47" "1.790.15"
$ws.Range("E47").Value = "  -0.96%  "

# Row 48 - RenderToken
$ws.Range("E48").Value = "  -2.06%  "

# Row 49 - Quant
Set-TextValue "D49" "88.21"
$ws.Range("E49").Value = "  -0.78%  "

# Row 50 - BabyDogeCoin
$ws.Range("D50").Value = "0.0₆0107"
$ws.Range("E50").Value = "  -0.50%  "

# Row 51 - Algorand
$ws.Range("E51").Value = "  -3.53%  "
